$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.076394557952881
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.382917404174805
$ws.Range("D1").Value = 1.346486449241638
$ws.Range("E1").Value = 0.9659369587898254
